$d = $word.ActiveDocument
[void]$d.Content.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:wpsCustomData="http://www.wps.cn/officeDocument/2013/wpsCustomData" mc:Ignorable="w14 w15 wp14"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>SALES CONTRACT</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Contract No: {{ contract_no }}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Date: {{ date }}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>The Buyer: {{ buyer_name }}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Address: {{ buyer_address }}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>The Seller: ZHENBAO MACHINERY CO., LTD</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Address: [Please Input Your Company Address Here]</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>This contract is made by and between the Buyer and the Seller; whereby the Buyer agrees to buy and the Seller agrees to sell the under-mentioned goods subject to the terms and conditions as stipulated below:</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblW w:w="0" w:type="auto"/><w:tblCellSpacing w:w="15" w:type="dxa"/><w:tblInd w:w="15" w:type="dxa"/><w:tblBorders><w:top w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:left w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:bottom w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:right w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:insideH w:val="none" w:color="auto" w:sz="0" w:space="0"/><w:insideV w:val="none" w:color="auto" w:sz="0" w:space="0"/></w:tblBorders><w:shd w:val="clear"/><w:tblLayout w:type="fixed"/><w:tblCellMar><w:top w:w="0" w:type="dxa"/><w:left w:w="0" w:type="dxa"/><w:bottom w:w="0" w:type="dxa"/><w:right w:w="0" w:type="dxa"/></w:tblCellMar></w:tblPr><w:tblGrid><w:gridCol w:w="2629"/><w:gridCol w:w="2300"/><w:gridCol w:w="1987"/><w:gridCol w:w="789"/><w:gridCol w:w="1036"/></w:tblGrid><w:tr><w:tblPrEx><w:tblBorders><w:top w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:left w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:bottom w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:right w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:insideH w:val="none" w:color="auto" w:sz="0" w:space="0"/><w:insideV w:val="none" w:color="auto" w:sz="0" w:space="0"/></w:tblBorders><w:shd w:val="clear"/><w:tblCellMar><w:top w:w="0" w:type="dxa"/><w:left w:w="0" w:type="dxa"/><w:bottom w:w="0" w:type="dxa"/><w:right w:w="0" w:type="dxa"/></w:tblCellMar></w:tblPrEx><w:trPr><w:tblCellSpacing w:w="15" w:type="dxa"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2584" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:left w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:bottom w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:right w:val="single" w:color="auto" w:sz="6" w:space="0"/></w:tcBorders><w:shd w:val="clear"/><w:tcMar><w:top w:w="120" w:type="dxa"/><w:left w:w="180" w:type="dxa"/><w:bottom w:w="120" w:type="dxa"/><w:right w:w="180" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>No.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2270" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:left w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:bottom w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:right w:val="single" w:color="auto" w:sz="6" w:space="0"/></w:tcBorders><w:shd w:val="clear"/><w:tcMar><w:top w:w="120" w:type="dxa"/><w:left w:w="180" w:type="dxa"/><w:bottom w:w="120" w:type="dxa"/><w:right w:w="180" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>Description</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1957" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:left w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:bottom w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:right w:val="single" w:color="auto" w:sz="6" w:space="0"/></w:tcBorders><w:shd w:val="clear"/><w:tcMar><w:top w:w="120" w:type="dxa"/><w:left w:w="180" w:type="dxa"/><w:bottom w:w="120" w:type="dxa"/><w:right w:w="180" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>Qty</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="759" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:left w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:bottom w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:right w:val="single" w:color="auto" w:sz="6" w:space="0"/></w:tcBorders><w:shd w:val="clear"/><w:tcMar><w:top w:w="120" w:type="dxa"/><w:left w:w="180" w:type="dxa"/><w:bottom w:w="120" w:type="dxa"/><w:right w:w="180" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>Unit Price</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="991" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:left w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:bottom w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:right w:val="single" w:color="auto" w:sz="6" w:space="0"/></w:tcBorders><w:shd w:val="clear"/><w:tcMar><w:top w:w="120" w:type="dxa"/><w:left w:w="180" w:type="dxa"/><w:bottom w:w="120" w:type="dxa"/><w:right w:w="180" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>Total Amount</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tblPrEx><w:tblBorders><w:top w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:left w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:bottom w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:right w:val="none" w:color="1F1F1F" w:sz="0" w:space="0"/><w:insideH w:val="none" w:color="auto" w:sz="0" w:space="0"/><w:insideV w:val="none" w:color="auto" w:sz="0" w:space="0"/></w:tblBorders><w:shd w:val="clear"/><w:tblCellMar><w:top w:w="0" w:type="dxa"/><w:left w:w="0" w:type="dxa"/><w:bottom w:w="0" w:type="dxa"/><w:right w:w="0" w:type="dxa"/></w:tblCellMar></w:tblPrEx><w:trPr><w:tblCellSpacing w:w="15" w:type="dxa"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2584" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:left w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:bottom w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:right w:val="single" w:color="auto" w:sz="6" w:space="0"/></w:tcBorders><w:shd w:val="clear"/><w:tcMar><w:top w:w="120" w:type="dxa"/><w:left w:w="180" w:type="dxa"/><w:bottom w:w="120" w:type="dxa"/><w:right w:w="180" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{%tr for item in items %}{{ item.no }}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2270" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:left w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:bottom w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:right w:val="single" w:color="auto" w:sz="6" w:space="0"/></w:tcBorders><w:shd w:val="clear"/><w:tcMar><w:top w:w="120" w:type="dxa"/><w:left w:w="180" w:type="dxa"/><w:bottom w:w="120" w:type="dxa"/><w:right w:w="180" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{{ item.desc_en }} {{ item.desc_cn }}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1957" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:left w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:bottom w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:right w:val="single" w:color="auto" w:sz="6" w:space="0"/></w:tcBorders><w:shd w:val="clear"/><w:tcMar><w:top w:w="120" w:type="dxa"/><w:left w:w="180" w:type="dxa"/><w:bottom w:w="120" w:type="dxa"/><w:right w:w="180" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{{ item.qty }} {{ item.unit }}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="759" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:left w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:bottom w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:right w:val="single" w:color="auto" w:sz="6" w:space="0"/></w:tcBorders><w:shd w:val="clear"/><w:tcMar><w:top w:w="120" w:type="dxa"/><w:left w:w="180" w:type="dxa"/><w:bottom w:w="120" w:type="dxa"/><w:right w:w="180" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{{ item.price }}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="991" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:left w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:bottom w:val="single" w:color="auto" w:sz="6" w:space="0"/><w:right w:val="single" w:color="auto" w:sz="6" w:space="0"/></w:tcBorders><w:shd w:val="clear"/><w:tcMar><w:top w:w="120" w:type="dxa"/><w:left w:w="180" w:type="dxa"/><w:bottom w:w="120" w:type="dxa"/><w:right w:w="180" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{{ item.total }}{%tr endtr %}</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Total Amount: USD {{ total_amount }}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>1. Shipping Method: {{ shipping_method }}</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2. Payment Terms: {{ payment_terms }}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>3. Lead Time: {{ lead_time }}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>4. Packing: Export Standard Packing.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>5. Insurance: To be covered by the Buyer.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>The Seller:                                     The Buyer:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>ZHENBAO MACHINERY CO., LTD                      {{ buyer_name }}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>__________________________                      __________________________</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Signature &amp; Stamp                               Signature &amp; Stamp</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
